$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5503
$ws.Range("H3").Value = 8347
$ws.Range("J3").Value = 5863
$ws.Range("I4").Value = 1770
$ws.Range("J4").Value = 1277
$ws.Range("J5").Value = 453
$ws.Range("J6").Value = 7358
$ws.Range("H7").Value = 26012
$ws.Range("I7").Value = 26224
$ws.Range("J7").Value = 20454

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 177
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 356
$ws.Range("J6").Value = 432
$ws.Range("J7").Value = 1283

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 107
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J4").Value = 21
$ws.Range("J7").Value = 630

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 127
$ws.Range("J7").Value = 321

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 163
$ws.Range("J7").Value = 592
$ws.Range("J8").Value = 1283
$ws.Range("J10").Value = 137
$ws.Range("J15").Value = 223
$ws.Range("J19").Value = 596
$ws.Range("J23").Value = 195
$ws.Range("J25").Value = 103
$ws.Range("J27").Value = 122
$ws.Range("J34").Value = 100
$ws.Range("J36").Value = 284
$ws.Range("J37").Value = 630
$ws.Range("J42").Value = 848
$ws.Range("J44").Value = 158
$ws.Range("H46").Value = 60
$ws.Range("J47").Value = 157
$ws.Range("J51").Value = 252
$ws.Range("J53").Value = 280
$ws.Range("J54").Value = 393
$ws.Range("J55").Value = 278
$ws.Range("J60").Value = 126
$ws.Range("J67").Value = 777
$ws.Range("J72").Value = 82
$ws.Range("J73").Value = 198
$ws.Range("J76").Value = 296
$ws.Range("J78").Value = 251
$ws.Range("J79").Value = 585
$ws.Range("J84").Value = 176
$ws.Range("J86").Value = 123
$ws.Range("J89").Value = 266
$ws.Range("J90").Value = 223
$ws.Range("J95").Value = 307
$ws.Range("J96").Value = 243
$ws.Range("I98").Value = 185
$ws.Range("J98").Value = 143
$ws.Range("J99").Value = 321
$ws.Range("H101").Value = 26012
$ws.Range("I101").Value = 26224
$ws.Range("J101").Value = 20454

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 299
$ws.Range("J5").Value = 22
$ws.Range("J6").Value = 203
$ws.Range("J7").Value = 777

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J3").Value = 57
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 187
$ws.Range("J7").Value = 393

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 150
$ws.Range("J5").Value = 23
$ws.Range("J6").Value = 217
$ws.Range("J7").Value = 596

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 188
$ws.Range("J3").Value = 171
$ws.Range("J4").Value = 39
$ws.Range("J6").Value = 433
$ws.Range("J7").Value = 848

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 27
$ws.Range("J7").Value = 251

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J4").Value = 11
$ws.Range("J6").Value = 137
$ws.Range("J7").Value = 278

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("H3").Value = 12
$ws.Range("H7").Value = 60

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 243

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 207
$ws.Range("J7").Value = 585

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 91
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 185
$ws.Range("J6").Value = 187
$ws.Range("J7").Value = 592

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 100

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 22
$ws.Range("I4").Value = 11
$ws.Range("J6").Value = 88
$ws.Range("I7").Value = 185
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 266

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 66
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 251

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 82
